$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the date string shared by B1:B8, B10, B11 from "19/07/2023"/"17/06/2023"
# to "05/08/2023" (keeping B9 on the older "17/06/2023" date), while preserving
# each cell's existing (quote-prefixed text) style.
foreach ($r in 1..8) {
    $ws.Cells.Item($r, 2).Value = "'05/08/2023"
}
$ws.Cells.Item(10, 2).Value = "'05/08/2023"
$ws.Cells.Item(11, 2).Value = "'05/08/2023"

# Move the active selection to H9 (search bar / filter UI area).
$ws.Range("H9").Select()

# Configure page setup: A4 paper, portrait orientation.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
